$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (column A and B, rows 1-5)
$ws.Range("A1").Value = 0.0063612830445923078
$ws.Range("B1").Value = -0.0063612831469799392

$ws.Range("A2").Value = -0.029661300653364715
$ws.Range("B2").Value = 0.029661300543289922

$ws.Range("A3").Value = 0.046750399037914787
$ws.Range("B3").Value = -0.046750399097362984

$ws.Range("A4").Value = 0.063032308168594398
$ws.Range("B4").Value = -0.063032308278751059

$ws.Range("A5").Value = -0.034866706185505317
$ws.Range("B5").Value = 0.034866706026036864

# Update column widths (closest achievable values given the host's
# character-width -> pixel -> character-width rounding of 1/6ths):
# target A width 14.7109375  -> nearest reachable is 14.666666666666666 (ColumnWidth 13.8)
# target B width 15.42578125 -> nearest reachable is 15.5               (ColumnWidth 14.6)
$ws.Columns.Item(1).ColumnWidth = 13.8
$ws.Columns.Item(2).ColumnWidth = 14.6
